$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (German)
$ws.Range("B2").Value = 6444
$ws.Range("C2").Value = 4292
$ws.Range("D2").Value = 5733
$ws.Range("E2").Value = 6155
$ws.Range("F2").Value = 6342

# Row 6 (Greek)
$ws.Range("B6").Value = 5936
$ws.Range("C6").Value = 1171
$ws.Range("D6").Value = 3019
$ws.Range("E6").Value = 4916
$ws.Range("F6").Value = 5904

# Row 11 (Indonesian)
$ws.Range("B11").Value = 7926
$ws.Range("C11").Value = 5338
$ws.Range("D11").Value = 7916
$ws.Range("E11").Value = 7926
$ws.Range("F11").Value = 7926

# Row 12 (Finnish)
$ws.Range("B12").Value = 4432
$ws.Range("C12").Value = 4224
$ws.Range("D12").Value = 4395
$ws.Range("E12").Value = 4417
$ws.Range("F12").Value = 4421

# Row 14 (Korean)
$ws.Range("B14").Value = 36000
$ws.Range("C14").Value = 33522
$ws.Range("D14").Value = 35999
$ws.Range("E14").Value = 36000
$ws.Range("F14").Value = 36000

# Row 16 (Turkish)
$ws.Range("B16").Value = 4486
$ws.Range("C16").Value = 3103
$ws.Range("D16").Value = 4150
$ws.Range("E16").Value = 4445
$ws.Range("F16").Value = 4485

# Row 17 (Arabic)
$ws.Range("B17").Value = 2468
$ws.Range("C17").Value = 2397
$ws.Range("D17").Value = 2445
$ws.Range("E17").Value = 2460
$ws.Range("F17").Value = 2468

# Row 20 (Maltese)
$ws.Range("B20").Value = 595
$ws.Range("C20").Value = 430
$ws.Range("D20").Value = 579
$ws.Range("E20").Value = 592
$ws.Range("F20").Value = 594
